$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45988, 11577.6521601806, 10857.5505923417, 17691.4,   7037.5879275933,  8.48910499728997),
    @(45989, 11475.3170215859, 10206.7041866431, 9003.4,    6791.6782700679,  333.124269029623),
    @(45990, 4528.62839486405, 7174.07284542601, 9003.4,    6144.13394854895, 179.783616415623),
    @(45991, 4146.45002731865, 6956.59577426702, 9003.4,    5858.80671675908, 158.833437126088),
    @(45992, 11227.111966293,  10365.7396492797, 11040.26,  7119.82572049075, 268.5543904071),
    @(45993, 11775.3234136769, 10670.1463912902, 11040.26,  7488.55906966544, 296.601894206484),
    @(45994, 11775.3234136769, 10401.2906072152, 11040.26,  7488.55906966544, 285.399569870025),
    @(45995, 11775.3234136769, 10340.1587293136, 11040.26,  7488.55906966544, 282.852408290793),
    @(45996, 11775.3234136769, 9796.50596114354, 11040.26,  7488.55906966544, 260.20020961704),
    @(45997, 4605.75571869631, 6804.00481553918, 11040.26,  7098.38188713283, 119.255279278),
    @(45998, 4450.60726673227, 6913.73962985427, 11040.26,  7083.31456051356, 123.199757931993),
    @(45999, 4752.2157406543,  7357.95324486679, 11040.26,  7517.03745133792, 159.780445675196),
    @(46000, 12413.0231244078, 11529.4647511327, 11040.26,  7940.95011014115, 351.256452553078),
    @(46001, 12413.0231244078, 11709.3325287338, 11040.26,  7940.95011014115, 358.750943286457)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
